# Pre-process Shir data for decoding analysis
# - Apply wrap-text style + taller header row for column C ("included_manual")
# - Append new experiment rows (bat 2382 and bat 0194) to the inclusion list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: make it taller and wrap the "included_manual" header text
$ws.Rows.Item(1).RowHeight = 30
$ws.Range("C1").WrapText = $true

# New data rows 151-209 (bat_num 2382 then bat_num 194)
$ws.Range("A151").Value = "b2382_d190623"
$ws.Range("B151").Value = 2382
$ws.Range("C151").Value = 1

$ws.Range("A152").Value = "b2382_d190624"
$ws.Range("B152").Value = 2382
$ws.Range("C152").Value = 0.5
$ws.Range("D152").Value = "not sure how bad is this over-representation"

$ws.Range("A153").Value = "b2382_d190625"
$ws.Range("B153").Value = 2382
$ws.Range("C153").Value = 0.5
$ws.Range("D153").Value = "not sure how bad is this over-representation"

$ws.Range("A154").Value = "b2382_d190627"
$ws.Range("B154").Value = 2382
$ws.Range("C154").Value = 1

$ws.Range("A155").Value = "b2382_d190628"
$ws.Range("B155").Value = 2382
$ws.Range("C155").Value = 0

$ws.Range("A156").Value = "b2382_d190630"
$ws.Range("B156").Value = 2382
$ws.Range("C156").Value = 0

$ws.Range("A157").Value = "b2382_d190701"
$ws.Range("B157").Value = 2382
$ws.Range("C157").Value = 0

$ws.Range("A158").Value = "b2382_d190703"
$ws.Range("B158").Value = 2382
$ws.Range("C158").Value = 0

$ws.Range("A159").Value = "b2382_d190707"
$ws.Range("B159").Value = 2382
$ws.Range("C159").Value = 0

$ws.Range("A160").Value = "b2382_d190708"
$ws.Range("B160").Value = 2382
$ws.Range("C160").Value = 1

$ws.Range("A161").Value = "b2382_d190709"
$ws.Range("B161").Value = 2382
$ws.Range("C161").Value = 0

$ws.Range("A162").Value = "b2382_d190712"
$ws.Range("B162").Value = 2382
$ws.Range("C162").Value = 1

$ws.Range("A163").Value = "b2382_d190714"
$ws.Range("B163").Value = 2382
$ws.Range("C163").Value = 1

$ws.Range("A164").Value = "b2382_d190715"
$ws.Range("B164").Value = 2382
$ws.Range("C164").Value = 1

$ws.Range("A165").Value = "b2382_d190716"
$ws.Range("B165").Value = 2382
$ws.Range("C165").Value = 1

$ws.Range("A166").Value = "b2382_d190718"
$ws.Range("B166").Value = 2382
$ws.Range("C166").Value = 1

$ws.Range("A167").Value = "b2382_d190721"
$ws.Range("B167").Value = 2382
$ws.Range("C167").Value = 1

$ws.Range("A168").Value = "b2382_d190722"
$ws.Range("B168").Value = 2382
$ws.Range("C168").Value = 1

$ws.Range("A169").Value = "b2382_d190724"
$ws.Range("B169").Value = 2382
$ws.Range("C169").Value = 0.5
$ws.Range("D169").Value = "over-representation only at the very last bin"

$ws.Range("A170").Value = "b2382_d190725"
$ws.Range("B170").Value = 2382
$ws.Range("C170").Value = 1

$ws.Range("A171").Value = "b2382_d190728"
$ws.Range("B171").Value = 2382
$ws.Range("C171").Value = 1

$ws.Range("A172").Value = "b2382_d190729"
$ws.Range("B172").Value = 2382
$ws.Range("C172").Value = 1

$ws.Range("A173").Value = "b2382_d190730"
$ws.Range("B173").Value = 2382
$ws.Range("C173").Value = 1

$ws.Range("A174").Value = "b2382_d190731"
$ws.Range("B174").Value = 2382
$ws.Range("C174").Value = 1

$ws.Range("A175").Value = "b2382_d190801"
$ws.Range("B175").Value = 2382
$ws.Range("C175").Value = 1

$ws.Range("A176").Value = "b2382_d190804"
$ws.Range("B176").Value = 2382
$ws.Range("C176").Value = 0.5
$ws.Range("D176").Value = "not sure how bad is this over-representation"

$ws.Range("A177").Value = "b2382_d190805"
$ws.Range("B177").Value = 2382
$ws.Range("C177").Value = 1

$ws.Range("A178").Value = "b2382_d190807"
$ws.Range("B178").Value = 2382
$ws.Range("C178").Value = 0.5
$ws.Range("D178").Value = "not sure how bad is this over-representation"

$ws.Range("A179").Value = "b2382_d190808"
$ws.Range("B179").Value = 2382
$ws.Range("C179").Value = 1

$ws.Range("A180").Value = "b2382_d190811"
$ws.Range("B180").Value = 2382
$ws.Range("C180").Value = 1

$ws.Range("A181").Value = "b2382_d190812"
$ws.Range("B181").Value = 2382
$ws.Range("C181").Value = 1

$ws.Range("A182").Value = "b2382_d190813"
$ws.Range("B182").Value = 2382
$ws.Range("C182").Value = 1

$ws.Range("A183").Value = "b2382_d190814"
$ws.Range("B183").Value = 2382
$ws.Range("C183").Value = 0

$ws.Range("A184").Value = "b0194_d180429"
$ws.Range("B184").Value = 194
$ws.Range("D184").Value = "need to process this day"

$ws.Range("A185").Value = "b0194_d180501"
$ws.Range("B185").Value = 194
$ws.Range("C185").Value = 0

$ws.Range("A186").Value = "b0194_d180502"
$ws.Range("B186").Value = 194
$ws.Range("C186").Value = 0

$ws.Range("A187").Value = "b0194_d180503"
$ws.Range("B187").Value = 194
$ws.Range("C187").Value = 1

$ws.Range("A188").Value = "b0194_d180505"
$ws.Range("B188").Value = 194
$ws.Range("C188").Value = 1

$ws.Range("A189").Value = "b0194_d180507"
$ws.Range("B189").Value = 194
$ws.Range("C189").Value = 1

$ws.Range("A190").Value = "b0194_d180508"
$ws.Range("B190").Value = 194
$ws.Range("C190").Value = 0.5
$ws.Range("D190").Value = "not sure how bad is this over-representation"

$ws.Range("A191").Value = "b0194_d180509"
$ws.Range("B191").Value = 194
$ws.Range("C191").Value = 0.5
$ws.Range("D191").Value = "over-representation only at the very first bin"

$ws.Range("A192").Value = "b0194_d180510"
$ws.Range("B192").Value = 194
$ws.Range("C192").Value = 1

$ws.Range("A193").Value = "b0194_d180513"
$ws.Range("B193").Value = 194
$ws.Range("C193").Value = 1

$ws.Range("A194").Value = "b0194_d180514"
$ws.Range("B194").Value = 194
$ws.Range("C194").Value = 1

$ws.Range("A195").Value = "b0194_d180515"
$ws.Range("B195").Value = 194
$ws.Range("C195").Value = 0.5
$ws.Range("D195").Value = "not sure how bad is this over-representation"

$ws.Range("A196").Value = "b0194_d180516"
$ws.Range("B196").Value = 194
$ws.Range("C196").Value = 0.5
$ws.Range("D196").Value = "not sure how bad is this over-representation"

$ws.Range("A197").Value = "b0194_d180520"
$ws.Range("B197").Value = 194
$ws.Range("C197").Value = 0.5
$ws.Range("D197").Value = "not sure how bad is this over-representation"

$ws.Range("A198").Value = "b0194_d180521"
$ws.Range("B198").Value = 194
$ws.Range("C198").Value = 0.5

$ws.Range("A199").Value = "b0194_d180522"
$ws.Range("B199").Value = 194
$ws.Range("C199").Value = 0

$ws.Range("A200").Value = "b0194_d180528"
$ws.Range("B200").Value = 194
$ws.Range("C200").Value = 0

$ws.Range("A201").Value = "b0194_d180531"
$ws.Range("B201").Value = 194
$ws.Range("C201").Value = 0

$ws.Range("A202").Value = "b0194_d180604"
$ws.Range("B202").Value = 194
$ws.Range("C202").Value = 1

$ws.Range("A203").Value = "b0194_d180605"
$ws.Range("B203").Value = 194
$ws.Range("C203").Value = 0.5
$ws.Range("D203").Value = "over-representation only at the very first bin"

$ws.Range("A204").Value = "b0194_d180606"
$ws.Range("B204").Value = 194
$ws.Range("C204").Value = 0

$ws.Range("A205").Value = "b0194_d180607"
$ws.Range("B205").Value = 194
$ws.Range("C205").Value = 0

$ws.Range("A206").Value = "b0194_d180610"
$ws.Range("B206").Value = 194
$ws.Range("C206").Value = 0

$ws.Range("A207").Value = "b0194_d180611"
$ws.Range("B207").Value = 194
$ws.Range("C207").Value = 0

$ws.Range("A208").Value = "b0194_d180612"
$ws.Range("B208").Value = 194
$ws.Range("C208").Value = 0

$ws.Range("A209").Value = "b0194_d180614"
$ws.Range("B209").Value = 194
$ws.Range("C209").Value = 0

# Reset the active selection to A1 (original "C150" selection is stale after the edit)
$ws.Range("A1").Select() | Out-Null
